$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.033.92"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.171.83"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.39"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.19"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +33.30%  "
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.169.16"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.769"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +14.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.203"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.75"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.93%  "
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "35.10"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +8.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.825.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.758.28"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.189.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +14.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.71"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +9.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +12.70%  "
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +16.70%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +11.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.43"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.343.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.95%  "
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.31"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.31"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +24.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.973"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -20.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.196"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +43.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "526.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.94"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.144"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0904"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +30.24%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.422"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +17.38%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +9.45%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.709"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +19.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "150.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.32%  "
$ws.Range("E49").Value = "  +9.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.38"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +13.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.09%  "
